$wb = $excel.ActiveWorkbook

# Sheet "addListItem": A2 holds the user name "UsertwelveH" -> bump to "UsertwelveI"
$wsAddListItem = $wb.Worksheets.Item("addListItem")
$wsAddListItem.Range("A2").Value = "UsertwelveI"

# Sheet "createUser": A2 holds the numeric test user id 1041 -> bump to 1043
$wsCreateUser = $wb.Worksheets.Item("createUser")
$wsCreateUser.Range("A2").Value = 1043

$wb.Save()
